$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-26T07:03:15.210224+00:00"
$ws.Range("K3").Value = "2025-11-26T07:03:17.528074+00:00"
$ws.Range("K4").Value = "2025-11-26T07:03:17.528113+00:00"
$ws.Range("K5").Value = "2025-11-26T07:03:20.297265+00:00"
$ws.Range("K6").Value = "2025-11-26T07:03:22.674197+00:00"
$ws.Range("K7").Value = "2025-11-26T07:03:25.519572+00:00"
$ws.Range("K8").Value = "2025-11-26T07:03:25.519599+00:00"
$ws.Range("K9").Value = "2025-11-26T07:03:25.519618+00:00"
$ws.Range("K10").Value = "2025-11-26T07:03:28.239812+00:00"
$ws.Range("K11").Value = "2025-11-26T07:03:30.596197+00:00"
$ws.Range("K12").Value = "2025-11-26T07:03:32.936403+00:00"
$ws.Range("K13").Value = "2025-11-26T07:03:35.296026+00:00"
$ws.Range("K14").Value = "2025-11-26T07:03:38.055390+00:00"
$ws.Range("K15").Value = "2025-11-26T07:03:42.994691+00:00"
$ws.Range("K16").Value = "2025-11-26T07:03:42.994724+00:00"
$ws.Range("K17").Value = "2025-11-26T07:03:42.994744+00:00"
$ws.Range("K18").Value = "2025-11-26T07:03:45.801991+00:00"
$ws.Range("K19").Value = "2025-11-26T07:03:45.802031+00:00"
$ws.Range("K20").Value = "2025-11-26T07:03:45.802050+00:00"
$ws.Range("K21").Value = "2025-11-26T07:03:45.802066+00:00"
$ws.Range("K22").Value = "2025-11-26T07:03:48.063525+00:00"
$ws.Range("K23").Value = "2025-11-26T07:03:48.063554+00:00"
$ws.Range("K24").Value = "2025-11-26T07:03:50.316307+00:00"
$ws.Range("K25").Value = "2025-11-26T07:03:50.316336+00:00"
$ws.Range("K26").Value = "2025-11-26T07:03:50.316353+00:00"
$ws.Range("K27").Value = "2025-11-26T07:03:50.316370+00:00"
$ws.Range("K28").Value = "2025-11-26T07:03:52.601170+00:00"
$ws.Range("K29").Value = "2025-11-26T07:03:55.117445+00:00"
$ws.Range("K30").Value = "2025-11-26T07:03:57.918244+00:00"
$ws.Range("K31").Value = "2025-11-26T07:03:57.918274+00:00"
$ws.Range("K32").Value = "2025-11-26T07:04:03.048924+00:00"
$ws.Range("K33").Value = "2025-11-26T07:04:05.928816+00:00"
$ws.Range("K34").Value = "2025-11-26T07:04:05.928843+00:00"
